# [Outlook] Map sessionData snippets
# Adds 7 new rows to the "Snippets" table describing the Outlook
# SessionData APIs (getAllSessionData / setSessionData / getSessionData /
# removeSessionData / clearSessionData), mirroring the existing
# Class / Member Name / Member ID / SnippetId / MethodName layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 fresh rows right after the current last data row (244) so the
# new rows inherit the same cell style ("s=1") as the rest of the table,
# exactly like typing new rows under an Excel Table does.
$ws.Rows("245:251").Insert()

# New rows to append, in table-column order:
#   Class, Member Name, Member ID (methods only), SnippetId, MethodName
$newRows = @(
    @("AppointmentCompose", "sessionData", $null, "outlook-session-data-apis", "getAllSessionData"),
    @("MessageCompose",     "sessionData", $null, "outlook-session-data-apis", "getAllSessionData"),
    @("SessionData", "setAsync",    1, "outlook-session-data-apis", "setSessionData"),
    @("SessionData", "getAsync",    1, "outlook-session-data-apis", "getSessionData"),
    @("SessionData", "getAllAsync", 1, "outlook-session-data-apis", "getAllSessionData"),
    @("SessionData", "removeAsync", 1, "outlook-session-data-apis", "removeSessionData"),
    @("SessionData", "clearAsync",  1, "outlook-session-data-apis", "clearSessionData")
)

$startRow = 245
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    if ($null -ne $vals[2]) {
        $ws.Cells.Item($r, 3).Value = $vals[2]
    }
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}

# Grow the "Snippets" table (and its AutoFilter) to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E251"))

# Match the saved view: selection on the new last cell.
$ws.Range("E251").Select()
